{"js": "// \"Version 2.\" -> \"Version 1.\"\n//\n// Applied as three small, run-local edits (mirroring how Word itself would\n// coalesce/split runs for an in-place retype) rather than one big\n// search-and-replace, so the surrounding run/bookmark structure is\n// preserved:\n//   1. Retype \"Version\" over the \"Versi\"/\"on\" run split -> merges into a\n//      single run.\n//   2. Delete the now-redundant trailing \".\" run (its period slides into\n//      the \" 1.\" run created in step 3).\n//   3. Replace the lone digit \"2\" with \"1.\" in place, inside the existing\n//      \" 2\" run.\n\nconst body = context.document.body;\n\n// 1. Merge \"Versi\" + \"on\" into a single \"Version\" run.\nconst versionHits = body.search(\"Version\", { matchCase: true });\nawait context.sync();\nif (versionHits.items.length > 0) {\n  versionHits.items[0].insertText(\"Version\", \"Replace\");\n  await context.sync();\n}\n\n// 2. Remove the trailing standalone \".\" run.\nconst periodHits = body.search(\".\", { matchCase: true });\nawait context.sync();\nif (periodHits.items.length > 0) {\n  periodHits.items[0].delete();\n  await context.sync();\n}\n\n// 3. Turn the version number \"2\" into \"1.\" in place.\nconst digitHits = body.search(\"2\", { matchCase: true });\nawait context.sync();\nif (digitHits.items.length > 0) {\n  digitHits.items[0].insertText(\"1.\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# \"Version 2.\" -> \"Version 1.\"\n#\n# Applied as small, run-local Find/Replace edits (mirroring how Word itself\n# would coalesce/split runs for an in-place retype) rather than one big\n# search-and-replace, so the surrounding run/bookmark structure is\n# preserved:\n#   1. Retype \"Version\" over the \"Versi\"/\"on\" run split -> merges into a\n#      single run. (A same-value assignment is a no-op in this host, so we\n#      round-trip through a distinct placeholder to force the coalesce.)\n#   2. Delete the now-redundant trailing \".\" run (its period slides into\n#      the \" 1.\" run created in step 3).\n#   3. Replace the lone digit \"2\" with \"1.\" in place, inside the existing\n#      \" 2\" run.\n\n$d = $word.ActiveDocument\n\n# 1. Merge \"Versi\" + \"on\" into a single \"Version\" run.\n$find = $d.Content.Find\n$find.Text = \"Version\"\nif ($find.Execute()) {\n  $find.Parent.Text = \"Version#TMP#\"\n}\n$find2 = $d.Content.Find\n$find2.Text = \"Version#TMP#\"\nif ($find2.Execute()) {\n  $find2.Parent.Text = \"Version\"\n}\n\n# 2. Remove the trailing standalone \".\" run.\n$find3 = $d.Content.Find\n$find3.Text = \".\"\nif ($find3.Execute()) {\n  $find3.Parent.Delete()\n}\n\n# 3. Turn the version number \"2\" into \"1.\" in place.\n$find4 = $d.Content.Find\n$find4.Text = \"2\"\nif ($find4.Execute()) {\n  $find4.Parent.Text = \"1.\"\n}\n"}
